$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 292, shifting existing rows 292-318 down to 293-319.
$ws.Rows(292).Insert()

# Populate the newly inserted row 292 with the new record's data.
$ws.Cells.Item(292, 1).Value = 7
$ws.Cells.Item(292, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(292, 3).Value = "Ñuble"
$ws.Cells.Item(292, 4).Value = (Get-Date -Year 2022 -Month 9 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(292, 5).Value = 16
$ws.Cells.Item(292, 6).Value = 100112008
$ws.Cells.Item(292, 7).Value = "Coliflor"
$ws.Cells.Item(292, 8).Value = "Sin especificar"
$ws.Cells.Item(292, 9).Value = "Primera"
$ws.Cells.Item(292, 10).Value = 400
$ws.Cells.Item(292, 11).Value = 1200
$ws.Cells.Item(292, 12).Value = 1300
$ws.Cells.Item(292, 13).Value = 1250
$ws.Cells.Item(292, 14).Value = "`$/unidad"
$ws.Cells.Item(292, 15).Value = "Región del Maule"
$ws.Cells.Item(292, 16).Value = 1250
$ws.Cells.Item(292, 17).Value = 1
$ws.Cells.Item(292, 18).Value = "Hortaliza"
